# Fruta / hortaliza, semanal
# Insert two new weekly data rows at row 858 (pushing the existing rows
# 858-916 down to 860-918), then populate the two new rows with the
# latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 858.
$ws.Range("A858:A859").EntireRow.Insert()

# --- Row 858 ---------------------------------------------------------
$ws.Cells.Item(858, 1).Value = 8
$ws.Cells.Item(858, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(858, 3).Value = "Coquimbo"
$ws.Cells.Item(858, 4).Value = 44585
$ws.Cells.Item(858, 5).Value = 4
$ws.Cells.Item(858, 6).Value = 100112020
$ws.Cells.Item(858, 7).Value = "Tomate"
$ws.Cells.Item(858, 8).Value = "Semiduro"
$ws.Cells.Item(858, 9).Value = "Primera"
$ws.Cells.Item(858, 10).Value = 1000
$ws.Cells.Item(858, 11).Value = 9000
$ws.Cells.Item(858, 12).Value = 10000
$ws.Cells.Item(858, 13).Value = 9500
$ws.Cells.Item(858, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(858, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(858, 16).Value = 528
$ws.Cells.Item(858, 17).Value = 18
$ws.Cells.Item(858, 18).Value = "Hortaliza"

# --- Row 859 ---------------------------------------------------------
$ws.Cells.Item(859, 1).Value = 8
$ws.Cells.Item(859, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(859, 3).Value = "Coquimbo"
$ws.Cells.Item(859, 4).Value = 44585
$ws.Cells.Item(859, 5).Value = 4
$ws.Cells.Item(859, 6).Value = 100112020
$ws.Cells.Item(859, 7).Value = "Tomate"
$ws.Cells.Item(859, 8).Value = "Semiduro"
$ws.Cells.Item(859, 9).Value = "Segunda"
$ws.Cells.Item(859, 10).Value = 400
$ws.Cells.Item(859, 11).Value = 6500
$ws.Cells.Item(859, 12).Value = 7000
$ws.Cells.Item(859, 13).Value = 6750
$ws.Cells.Item(859, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(859, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(859, 16).Value = 375
$ws.Cells.Item(859, 17).Value = 18
$ws.Cells.Item(859, 18).Value = "Hortaliza"
